$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final target values for rows 2-16 (columns A: Name, B: Position, C: Team)
$data = @(
  @(2,  "Chris Paul",         "PG",    "San Antonio Spurs"),
  @(3,  "Jalen Green",        "PG,SG", "Houston Rockets"),
  @(4,  "Russell Westbrook",  "PG",    "Denver Nuggets"),
  @(6,  "Pascal Siakam",      "SF,PF", "Indiana Pacers"),
  @(7,  "Deni Avdija",        "SF,PF", "Portland Trail Blazers"),
  @(8,  "Bogdan Bogdanovic",  "SG,SF", "Atlanta Hawks"),
  @(9,  "Jerami Grant",       "SF,PF", "Portland Trail Blazers"),
  @(10, "Rudy Gobert",        "C",     "Minnesota Timberwolves"),
  @(14, "Corey Kispert",      "SG,SF", "Washington Wizards"),
  @(15, "Jaylen Brown",       "SG,SF", "Boston Celtics"),
  @(16, "Klay Thompson",      "SG,SF", "Dallas Mavericks")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}
